# Updates cryptos list values (Price column D, Volume(1h) column E) for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the literal text into the cell without Excel re-typing it as a number
    # (e.g. "108.29" -> 108.29) or losing trailing zeros / precision, and without
    # leaving a lasting NumberFormat change on the cell once done.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '52.199.60'
$ws.Range('E2').Value = '  +0.44%  '

$ws.Range('D3').Value = '2.988.41'
$ws.Range('E3').Value = '  +0.60%  '

$ws.Range('E4').Value = '  -0.03%  '

Set-TextValue $ws.Range('D5') '353.88'
$ws.Range('E5').Value = '  -0.01%  '

Set-TextValue $ws.Range('D6') '108.29'
$ws.Range('E6').Value = '  -3.70%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  -0.02%  '

Set-TextValue $ws.Range('D9') '0.623'
$ws.Range('E9').Value = '  -1.30%  '

$ws.Range('E10').Value = '  -3.42%  '

$ws.Range('E11').Value = '  +2.50%  '

$ws.Range('E12').Value = '  -4.12%  '

$ws.Range('E13').Value = '  -3.21%  '

$ws.Range('D14').Value = '3.462.67'
$ws.Range('E14').Value = '  +0.57%  '

Set-TextValue $ws.Range('D15') '7.77'
$ws.Range('E15').Value = '  -2.55%  '

$ws.Range('D16').Value = '2.991.33'
$ws.Range('E16').Value = '  +0.10%  '

$ws.Range('E17').Value = '  +3.19%  '

$ws.Range('D18').Value = '52.213.50'
$ws.Range('E18').Value = '  +0.26%  '

Set-TextValue $ws.Range('D19') '3.50'
$ws.Range('E19').Value = '  +5.61%  '

Set-TextValue $ws.Range('D20') '7.55'
$ws.Range('E20').Value = '  -1.99%  '

Set-TextValue $ws.Range('D21') '13.65'
$ws.Range('E21').Value = '  -6.08%  '

$ws.Range('D22').Value = '0.0₃0975'
$ws.Range('E22').Value = '  -1.54%  '

Set-TextValue $ws.Range('D23') '69.56'
$ws.Range('E23').Value = '  -2.55%  '

Set-TextValue $ws.Range('D24') '264.30'
$ws.Range('E24').Value = '  -2.34%  '

$ws.Range('E25').Value = '  -2.48%  '

$ws.Range('E26').Value = '  +0.65%  '

Set-TextValue $ws.Range('D27') '26.85'
$ws.Range('E27').Value = '  -3.16%  '

Set-TextValue $ws.Range('D28') '7.50'
$ws.Range('E28').Value = '  -2.19%  '

$ws.Range('E29').Value = '  -0.08%  '

Set-TextValue $ws.Range('D30') '0.108'
$ws.Range('E30').Value = '  -2.27%  '

$ws.Range('E31').Value = '  -3.76%  '

Set-TextValue $ws.Range('D32') '6.37'
$ws.Range('E32').Value = '  +2.29%  '

Set-TextValue $ws.Range('D33') '36.65'
$ws.Range('E33').Value = '  -2.82%  '

$ws.Range('E34').Value = '  +9.31%  '

Set-TextValue $ws.Range('D35') '50.92'
$ws.Range('E35').Value = '  -3.71%  '

Set-TextValue $ws.Range('D36') '0.0447'
$ws.Range('E36').Value = '  -0.61%  '

$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('E38').Value = '  -7.01%  '

Set-TextValue $ws.Range('D39') '17.89'
$ws.Range('E39').Value = '  -5.91%  '

$ws.Range('E40').Value = '  -4.75%  '

$ws.Range('E41').Value = '  +0.76%  '

$ws.Range('E42').Value = '  -0.34%  '

Set-TextValue $ws.Range('D43') '22.78'
$ws.Range('E43').Value = '  -4.79%  '

Set-TextValue $ws.Range('D44') '122.31'
$ws.Range('E44').Value = '  +7.34%  '

$ws.Range('E45').Value = '  -1.42%  '

$ws.Range('D46').Value = '2.125.02'
$ws.Range('E46').Value = '  -2.54%  '

Set-TextValue $ws.Range('D47') '3.39'
$ws.Range('E47').Value = '  -4.50%  '

$ws.Range('E48').Value = '  -5.25%  '

Set-TextValue $ws.Range('D49') '0.248'
$ws.Range('E49').Value = '  +1.73%  '

Set-TextValue $ws.Range('D50') '0.0331'
$ws.Range('E50').Value = '  -3.23%  '

Set-TextValue $ws.Range('D51') '0.925'
$ws.Range('E51').Value = '  -1.44%  '
